# Working Hours.xlsx - "Minor changes in attendance servlet"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column H to fit the new, longer text entries
$ws.Columns.Item(8).ColumnWidth = 103.6

# Row 86: add afternoon hours and a note about validation work
$ws.Range("C86").Value = 4.5
$ws.Range("H86").Value = "Validation text in editAddStudent & editDelStudent Servlets and validations in addSubject & addFaculty Servlets "
$ws.Range("J86").Value = 4

# Row 87: add afternoon hours and a note about the attendance servlet
$ws.Range("C87").Value = 1
$ws.Range("H87").Value = "attendance servlet"

# Row 88: add afternoon hours and a note about the attendance servlet
$ws.Range("C88").Value = 5
$ws.Range("H88").Value = "attendance servlet"

# Extend the totals row formulas to cover the newly-filled rows
$ws.Range("E96").Formula = "=SUM(E3:E94)"
$ws.Range("J96").Formula = "=SUM(J3:J94)"

# Reflect the freshly-edited rows in the view state
$ws.Range("H88").Select()
